# Apply updated cryptocurrency price/volume figures scraped on
# Sat Mar  4 17:50:16 UTC 2023 (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.371.59"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "'1.567.12"
$ws.Range("D4").Value = "'1.002"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").Value = "'291.07"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("D7").Value = "'0.3782"
$ws.Range("E7").Value = "  +2.86%  "
$ws.Range("D8").Value = "'49.24"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.3401"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("D11").Value = "'1.139"
$ws.Range("E11").Value = "  -2.94%  "
$ws.Range("D12").Value = "'1.002"
$ws.Range("D13").Value = "'21.06"
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").Value = "'5.988"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").Value = "'6.917"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").Value = "'1.566.93"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").Value = "'0.00001134"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "'89.96"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").Value = "'0.06744"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").Value = "'16.62"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "'6.200"
$ws.Range("E22").Value = "  -1.18%  "
$ws.Range("D23").Value = "'11.95"
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("D24").Value = "'22.374.62"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").Value = "'2.401"
$ws.Range("E25").Value = "  +1.99%  "
$ws.Range("D26").Value = "'2.690"
$ws.Range("E26").Value = "  -7.26%  "
$ws.Range("D27").Value = "'20.11"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "'147.49"
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("D29").Value = "'5.029"
$ws.Range("E29").Value = "  +0.88%  "
$ws.Range("D30").Value = "'125.99"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").Value = "'1.739.60"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "'2.016"
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("D33").Value = "'6.082"
$ws.Range("E33").Value = "  -2.66%  "
$ws.Range("D34").Value = "'0.9942"
$ws.Range("E34").Value = "  -3.27%  "
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").Value = "'1.435"
$ws.Range("E36").Value = "  +10.14%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "'0.02510"
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("D39").Value = "'0.2292"
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("D40").Value = "'0.06486"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").Value = "'5.407"
$ws.Range("E41").Value = "  -2.46%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.6320"
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'11.32"
$ws.Range("E43").Value = "  -3.63%  "
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").Value = "'13.93"
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("D46").Value = "'3.800"
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("D47").Value = "'0.5937"
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("D48").Value = "'2.085"
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").Value = "'1.257"
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("D50").Value = "'124.69"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").Value = "'0.07323"
$ws.Range("E51").Value = "  +0.35%  "

# The apostrophe prefix above marks the cells with a "quote prefix" style;
# restore the column to its original (default) style now that the values
# are safely stored as text.
$ws.Range("D2:D51").Style = "Normal"

